$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 195; this shifts the existing rows 195-204
# down to 196-205 (and copies row formatting, e.g. the date style on D).
$ws.Rows("195:195").Insert()

# Populate the newly inserted row 195 with the new weekly price entry.
$ws.Cells.Item(195, 1).Value = 7
$ws.Cells.Item(195, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(195, 3).Value = "Ñuble"
$ws.Cells.Item(195, 4).Value = 44516
$ws.Cells.Item(195, 5).Value = 16
$ws.Cells.Item(195, 6).Value = 100114013
$ws.Cells.Item(195, 7).Value = "Zanahoria"
$ws.Cells.Item(195, 8).Value = "Sin especificar"
$ws.Cells.Item(195, 9).Value = "Primera"
$ws.Cells.Item(195, 10).Value = 80
$ws.Cells.Item(195, 11).Value = 8000
$ws.Cells.Item(195, 12).Value = 8500
$ws.Cells.Item(195, 13).Value = 8250
$ws.Cells.Item(195, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(195, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(195, 16).Value = 412
$ws.Cells.Item(195, 17).Value = 20
$ws.Cells.Item(195, 18).Value = "Hortaliza"
